$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2710.5557
$ws.Range("I100").Value = 2199.3333
$ws.Range("J100").Value = 2966.1667
$ws.Range("K100").Value = 2199.3333
$ws.Range("L100").Value = 2966.1667
$ws.Range("M100").Value = -1658.3333
$ws.Range("N100").Value = -4048.1667
$ws.Range("H112").Value = 6811
$ws.Range("I112").Value = 766.6667
$ws.Range("J112").Value = 7205.196
$ws.Range("K112").Value = 2300.0001
$ws.Range("L112").Value = 21615.588
$ws.Range("M112").Value = -1192.0001
$ws.Range("N112").Value = -23831.588
$ws.Range("H113").Value = 2234.5652
$ws.Range("I113").Value = 1927.7273
$ws.Range("J113").Value = 2515.8333
$ws.Range("K113").Value = 1927.7273
$ws.Range("L113").Value = 2515.8333
$ws.Range("M113").Value = 1326.2727
$ws.Range("N113").Value = -9023.8333
$ws.Range("H129").Value = 1253.8541
$ws.Range("I129").Value = 799
$ws.Range("J129").Value = 1273.6305
$ws.Range("K129").Value = 2397
$ws.Range("L129").Value = 3820.8915
$ws.Range("M129").Value = 2603
$ws.Range("N129").Value = -13820.8915
$ws.Range("H137").Value = 2875784.5
$ws.Range("I137").Value = 4903639
$ws.Range("J137").Value = 2990.5
$ws.Range("K137").Value = 14710917
$ws.Range("L137").Value = 8971.5
$ws.Range("M137").Value = -14708367
$ws.Range("N137").Value = -14071.5
$ws.Range("H138").Value = 4648.217
$ws.Range("I138").Value = 5268.5835
$ws.Range("J138").Value = 4493.125
$ws.Range("K138").Value = 15805.7505
$ws.Range("L138").Value = 13479.375
$ws.Range("M138").Value = -10665.7505
$ws.Range("N138").Value = -23759.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1668.0714
$ws.Range("I2").Value = 1759
$ws.Range("J2").Value = 1440.75
$ws.Range("K2").Value = 1759
$ws.Range("L2").Value = 1440.75
$ws.Range("M2").Value = -1646
$ws.Range("N2").Value = -1666.75
$ws.Range("H61").Value = 16669760
$ws.Range("I61").Value = 25643286
$ws.Range("J61").Value = 4642.857
$ws.Range("K61").Value = 25643286
$ws.Range("L61").Value = 4642.857
$ws.Range("M61").Value = -25643074
$ws.Range("N61").Value = -5066.857
$ws.Range("H76").Value = 69396
$ws.Range("J76").Value = 69396
$ws.Range("L76").Value = 69396
$ws.Range("N76").Value = -70072
$ws.Range("H79").Value = 69396
$ws.Range("J79").Value = 69396
$ws.Range("L79").Value = 69396
$ws.Range("N79").Value = -71736
$ws.Range("H92").Value = 48062.57
$ws.Range("J92").Value = 48062.57
$ws.Range("L92").Value = 48062.57
$ws.Range("N92").Value = -53054.57
$ws.Range("H103").Value = 69181
$ws.Range("J103").Value = 69181
$ws.Range("L103").Value = 69181
$ws.Range("N103").Value = -71525
$ws.Range("H110").Value = 150875
$ws.Range("I110").Value = 225850
$ws.Range("J110").Value = 75900
$ws.Range("K110").Value = 225850
$ws.Range("L110").Value = 75900
$ws.Range("M110").Value = -223805
$ws.Range("N110").Value = -79990
$ws.Range("H116").Value = 1668.0714
$ws.Range("I116").Value = 1759
$ws.Range("J116").Value = 1440.75
$ws.Range("K116").Value = 1759
$ws.Range("L116").Value = 1440.75
$ws.Range("M116").Value = 535
$ws.Range("N116").Value = -6028.75
$ws.Range("H122").Value = 78692.766
$ws.Range("I122").Value = 92635.63
$ws.Range("K122").Value = 277906.89
$ws.Range("M122").Value = -275456.89
$ws.Range("H136").Value = 16669760
$ws.Range("I136").Value = 25643286
$ws.Range("J136").Value = 4642.857
$ws.Range("K136").Value = 76929858
$ws.Range("L136").Value = 13928.571
$ws.Range("M136").Value = -76927308
$ws.Range("N136").Value = -19028.571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1668.0714
$ws.Range("I3").Value = 1759
$ws.Range("J3").Value = 1440.75
$ws.Range("K3").Value = 1759
$ws.Range("L3").Value = 1440.75
$ws.Range("M3").Value = -1645
$ws.Range("N3").Value = -1668.75
$ws.Range("H20").Value = 62502360
$ws.Range("I20").Value = 2443.3
$ws.Range("J20").Value = 166668880
$ws.Range("K20").Value = 2443.3
$ws.Range("L20").Value = 166668880
$ws.Range("M20").Value = -2196.3
$ws.Range("N20").Value = -166669374
$ws.Range("H140").Value = 51256
$ws.Range("J140").Value = 51256
$ws.Range("L140").Value = 51256
$ws.Range("N140").Value = -61616

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6506.8096
$ws.Range("I31").Value = 2472.7878
$ws.Range("J31").Value = 9117.058999999999
$ws.Range("K31").Value = 2472.7878
$ws.Range("L31").Value = 9117.058999999999
$ws.Range("M31").Value = -2177.7878
$ws.Range("N31").Value = -9707.058999999999
$ws.Range("H34").Value = 6506.8096
$ws.Range("I34").Value = 2472.7878
$ws.Range("J34").Value = 9117.058999999999
$ws.Range("K34").Value = 2472.7878
$ws.Range("L34").Value = 9117.058999999999
$ws.Range("M34").Value = -2270.7878
$ws.Range("N34").Value = -9521.058999999999
$ws.Range("H134").Value = 6761125.5
$ws.Range("I134").Value = 7147304
$ws.Range("K134").Value = 21441912
$ws.Range("M134").Value = -21439377

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5066.36
$ws.Range("I3").Value = 3333.45
$ws.Range("J3").Value = 11998
$ws.Range("K3").Value = 10000.35
$ws.Range("L3").Value = 35994
$ws.Range("M3").Value = -9888.349999999999
$ws.Range("N3").Value = -36218
$ws.Range("H6").Value = 333.63635
$ws.Range("I6").Value = 36.583332
$ws.Range("J6").Value = 690.1
$ws.Range("K6").Value = 109.749996
$ws.Range("L6").Value = 2070.3
$ws.Range("M6").Value = 3.250004000000004
$ws.Range("N6").Value = -2296.3
$ws.Range("H10").Value = 684.5
$ws.Range("I10").Value = 377.85715
$ws.Range("J10").Value = 1400
$ws.Range("K10").Value = 1133.57145
$ws.Range("L10").Value = 4200
$ws.Range("M10").Value = -994.5714499999999
$ws.Range("N10").Value = -4478
$ws.Range("H11").Value = 176711.3
$ws.Range("I11").Value = 333462.34
$ws.Range("J11").Value = 366.375
$ws.Range("K11").Value = 1000387.02
$ws.Range("L11").Value = 1099.125
$ws.Range("M11").Value = -1000247.02
$ws.Range("N11").Value = -1379.125
$ws.Range("H106").Value = 8514.5
$ws.Range("J106").Value = 8514.5
$ws.Range("L106").Value = 25543.5
$ws.Range("N106").Value = -27435.5
$ws.Range("H113").Value = 719.91174
$ws.Range("I113").Value = 713.05554
$ws.Range("J113").Value = 727.625
$ws.Range("K113").Value = 2139.16662
$ws.Range("L113").Value = 2182.875
$ws.Range("M113").Value = 30.83338000000003
$ws.Range("N113").Value = -6522.875
$ws.Range("H124").Value = 2578.077
$ws.Range("I124").Value = 460.5
$ws.Range("J124").Value = 2963.0908
$ws.Range("K124").Value = 1381.5
$ws.Range("L124").Value = 8889.2724
$ws.Range("M124").Value = 3528.5
$ws.Range("N124").Value = -18709.2724
$ws.Range("H129").Value = 979075.5600000001
$ws.Range("I129").Value = 557.3333
$ws.Range("J129").Value = 1213920
$ws.Range("K129").Value = 1671.9999
$ws.Range("L129").Value = 3641760
$ws.Range("M129").Value = 3328.0001
$ws.Range("N129").Value = -3651760
$ws.Range("H131").Value = 3514.5435
$ws.Range("I131").Value = 583.1667
$ws.Range("J131").Value = 4549.147
$ws.Range("K131").Value = 1749.5001
$ws.Range("L131").Value = 13647.441
$ws.Range("M131").Value = 3290.4999
$ws.Range("N131").Value = -23727.441

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2033.8235
$ws.Range("I122").Value = 2028.8462
$ws.Range("K122").Value = 6086.5386
$ws.Range("M122").Value = -3636.5386

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1567.76
$ws.Range("I68").Value = 1530
$ws.Range("J68").Value = 1592.9333
$ws.Range("K68").Value = 1530
$ws.Range("L68").Value = 1592.9333
$ws.Range("M68").Value = -781
$ws.Range("N68").Value = -3090.9333
$ws.Range("H70").Value = 44081.75
$ws.Range("J70").Value = 44081.75
$ws.Range("L70").Value = 44081.75
$ws.Range("N70").Value = -44621.75
$ws.Range("H71").Value = 1567.76
$ws.Range("I71").Value = 1530
$ws.Range("J71").Value = 1592.9333
$ws.Range("K71").Value = 7650
$ws.Range("L71").Value = 7964.666499999999
$ws.Range("M71").Value = -3906
$ws.Range("N71").Value = -15452.6665
$ws.Range("H73").Value = 44081.75
$ws.Range("J73").Value = 44081.75
$ws.Range("L73").Value = 44081.75
$ws.Range("N73").Value = -45953.75
$ws.Range("H74").Value = 23571.428
$ws.Range("I74").Value = 30000
$ws.Range("J74").Value = 22500
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 22500
$ws.Range("M74").Value = -29002
$ws.Range("N74").Value = -24496
$ws.Range("H77").Value = 23571.428
$ws.Range("I77").Value = 30000
$ws.Range("J77").Value = 22500
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 67500
$ws.Range("M77").Value = -85008
$ws.Range("N77").Value = -77484
$ws.Range("H82").Value = 1468.3636
$ws.Range("I82").Value = 1471.7142
$ws.Range("J82").Value = 1462.5
$ws.Range("K82").Value = 1471.7142
$ws.Range("L82").Value = 1462.5
$ws.Range("M82").Value = -1110.7142
$ws.Range("N82").Value = -2184.5
$ws.Range("H85").Value = 1468.3636
$ws.Range("I85").Value = 1471.7142
$ws.Range("J85").Value = 1462.5
$ws.Range("K85").Value = 1471.7142
$ws.Range("L85").Value = 1462.5
$ws.Range("M85").Value = -223.7141999999999
$ws.Range("N85").Value = -3958.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14073.333
$ws.Range("J54").Value = 14073.333
$ws.Range("L54").Value = 14073.333
$ws.Range("N54").Value = -15113.333
$ws.Range("H60").Value = 90000
$ws.Range("J60").Value = 90000
$ws.Range("L60").Value = 90000
$ws.Range("N60").Value = -91644
$ws.Range("H81").Value = 5163.8667
$ws.Range("I81").Value = 6583.875
$ws.Range("J81").Value = 3541
$ws.Range("K81").Value = 13167.75
$ws.Range("L81").Value = 7082
$ws.Range("M81").Value = -12106.75
$ws.Range("N81").Value = -9204
$ws.Range("H82").Value = 39980
$ws.Range("J82").Value = 39980
$ws.Range("L82").Value = 39980
$ws.Range("N82").Value = -40746
$ws.Range("H84").Value = 5163.8667
$ws.Range("I84").Value = 6583.875
$ws.Range("J84").Value = 3541
$ws.Range("K84").Value = 65838.75
$ws.Range("L84").Value = 35410
$ws.Range("M84").Value = -60534.75
$ws.Range("N84").Value = -46018
$ws.Range("H85").Value = 39980
$ws.Range("J85").Value = 39980
$ws.Range("L85").Value = 39980
$ws.Range("N85").Value = -42632
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
